$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2025-01-14 Tuesday" "2025-01-15 Wednesday"

Replace-Text "867×8=6936" "432×5=2160"
Replace-Text "222×4=888" "312×8=2496"
Replace-Text "364×5=1820" "240×3=720"
Replace-Text "598×4=2392" "374×2=748"
Replace-Text "850×9=7650" "700×3=2100"

Replace-Text "651×5=3255" "155×7=1085"
Replace-Text "603×7=4221" "635×2=1270"
Replace-Text "203×6=1218" "681×7=4767"
Replace-Text "623×5=3115" "381×2=762"
Replace-Text "336×6=2016" "273×8=2184"

Replace-Text "977×7=6839" "431×3=1293"
Replace-Text "593×5=2965" "333×9=2997"
Replace-Text "166×5=830" "958×2=1916"
Replace-Text "936×4=3744" "565×4=2260"
Replace-Text "768×5=3840" "835×7=5845"

Replace-Text "869×5=4345" "704×6=4224"
Replace-Text "118×9=1062" "229×5=1145"
Replace-Text "539×8=4312" "931×3=2793"
Replace-Text "988×4=3952" "751×3=2253"
Replace-Text "781×8=6248" "564×4=2256"

Replace-Text "279×9=2511" "109×2=218"
Replace-Text "766×5=3830" "865×5=4325"
Replace-Text "250×3=750" "588×3=1764"
Replace-Text "610×9=5490" "907×3=2721"
Replace-Text "501×8=4008" "788×4=3152"
